$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. "Cost Breakdowns" sheet: the "Other" equipment cost-breakdown rows
#    (11 = ISIC 27 share, 12 = ISIC 41T43 share) previously had hardcoded
#    values (1/1/1 and 0/0/0). They are now driven by formulas that reuse
#    the "Appliances" row's breakdown (rows 9/10), and the explanatory
#    note text in column G is updated accordingly.
# -----------------------------------------------------------------------
$wsCB = $wb.Worksheets.Item("Cost Breakdowns")

$wsCB.Range("C11").Formula = "=C9"
$wsCB.Range("D11").Formula = "=D9"
$wsCB.Range("E11").Formula = "=E9"

$wsCB.Range("C12").Formula = "=C10"
$wsCB.Range("D12").Formula = "=D10"
$wsCB.Range("E12").Formula = "=E10"

$wsCB.Range("G11").Value = "We use the same breakdown as for appliances for"
$wsCB.Range("G12").Value = "the other equipment category"

# -----------------------------------------------------------------------
# 2. The three "SoBCaICbIC-*" sheets shared an "ISIC 20T21" combined
#    column. It is now split into two separate columns, "ISIC 20" and
#    "ISIC 21" -- achieved by inserting a new column before the old
#    second ISIC column (column L) and relabeling the headers.
# -----------------------------------------------------------------------
$sheetNames = @(
    "SoBCaICbIC-urbanresidential",
    "SoBCaICbIC-ruralresidential",
    "SoBCaICbIC-commercial"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns("L").Insert()
    $ws.Range("K1").Value = "ISIC 20"
    $ws.Range("L1").Value = "ISIC 21"
}
